$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Classi" (sheet1.xml) -- add row 7
# ---------------------------------------------------------------------------
$wsClassi = $wb.Worksheets.Item("Classi")

$wsClassi.Cells.Item(7, 1).Value = "14-19 febbraio 2022"
$wsClassi.Cells.Item(7, 2).Value = 6134
$wsClassi.Cells.Item(7, 2).NumberFormat = "#,##0"
$wsClassi.Cells.Item(7, 3).Value = 8157
$wsClassi.Cells.Item(7, 3).NumberFormat = "#,##0"
$wsClassi.Cells.Item(7, 4).Value = 0.752
$wsClassi.Cells.Item(7, 4).NumberFormat = "0.0%"
$wsClassi.Cells.Item(7, 5).Value = 376215
$wsClassi.Cells.Item(7, 5).NumberFormat = "#,##0"
$wsClassi.Cells.Item(7, 6).Value = 282974
$wsClassi.Cells.Item(7, 6).NumberFormat = "#,##0"
$wsClassi.Cells.Item(7, 7).Value = 0.752
$wsClassi.Cells.Item(7, 7).NumberFormat = "0.0%"
$wsClassi.Cells.Item(7, 8).Value = 282485
$wsClassi.Cells.Item(7, 8).NumberFormat = "#,##0"
$wsClassi.Cells.Item(7, 9).Value = 26787
$wsClassi.Cells.Item(7, 9).NumberFormat = "#,##0"
$wsClassi.Cells.Item(7, 10).Value = 0.998
$wsClassi.Cells.Item(7, 10).NumberFormat = "0.0%"
$wsClassi.Cells.Item(7, 11).Value = 0.095
$wsClassi.Cells.Item(7, 11).NumberFormat = "0.0%"
$wsClassi.Cells.Item(7, 14).Value = 489
$wsClassi.Cells.Item(7, 15).Value = 0.2

# ---------------------------------------------------------------------------
# Sheet "Alunni in presenza" (sheet2.xml) -- add row 7
# ---------------------------------------------------------------------------
$wsAlunniPresenza = $wb.Worksheets.Item("Alunni in presenza")

$wsAlunniPresenza.Cells.Item(7, 1).Value = "14-19 febbraio 2022"
$wsAlunniPresenza.Cells.Item(7, 2).Value = 7388444
$wsAlunniPresenza.Cells.Item(7, 2).NumberFormat = "#,##0"
$wsAlunniPresenza.Cells.Item(7, 3).Value = 5536315
$wsAlunniPresenza.Cells.Item(7, 3).NumberFormat = "#,##0"
$wsAlunniPresenza.Cells.Item(7, 4).Value = 0.749
$wsAlunniPresenza.Cells.Item(7, 4).NumberFormat = "0.0%"
$wsAlunniPresenza.Cells.Item(7, 5).Value = 5311636
$wsAlunniPresenza.Cells.Item(7, 5).NumberFormat = "#,##0"
$wsAlunniPresenza.Cells.Item(7, 6).Value = 0.959
$wsAlunniPresenza.Cells.Item(7, 6).NumberFormat = "0.0%"

# ---------------------------------------------------------------------------
# Sheet "Alunni" (sheet3.xml) -- add rows 22-24
# ---------------------------------------------------------------------------
$wsAlunni = $wb.Worksheets.Item("Alunni")

$wsAlunni.Cells.Item(22, 1).Value = "14-19 febbraio 2022"
$wsAlunni.Cells.Item(22, 2).Value = "Infanzia"
$wsAlunni.Cells.Item(22, 3).Value = 632694
$wsAlunni.Cells.Item(22, 3).NumberFormat = "#,##0"
$wsAlunni.Cells.Item(22, 4).Value = 606919
$wsAlunni.Cells.Item(22, 4).NumberFormat = "#,##0"
$wsAlunni.Cells.Item(22, 5).Value = 25775
$wsAlunni.Cells.Item(22, 5).NumberFormat = "#,##0"
$wsAlunni.Cells.Item(22, 6).Value = 0.041
$wsAlunni.Cells.Item(22, 6).NumberFormat = "0.0%"

$wsAlunni.Cells.Item(23, 1).Value = "14-19 febbraio 2022"
$wsAlunni.Cells.Item(23, 2).Value = "Primaria"
$wsAlunni.Cells.Item(23, 3).Value = 1745692
$wsAlunni.Cells.Item(23, 3).NumberFormat = "#,##0"
$wsAlunni.Cells.Item(23, 4).Value = 1670216
$wsAlunni.Cells.Item(23, 4).NumberFormat = "#,##0"
$wsAlunni.Cells.Item(23, 5).Value = 75476
$wsAlunni.Cells.Item(23, 5).NumberFormat = "#,##0"
$wsAlunni.Cells.Item(23, 6).Value = 0.043
$wsAlunni.Cells.Item(23, 6).NumberFormat = "0.0%"

$wsAlunni.Cells.Item(24, 1).Value = "14-19 febbraio 2022"
$wsAlunni.Cells.Item(24, 2).Value = "Sec. 1° e 2° Grado"
$wsAlunni.Cells.Item(24, 3).Value = 3157929
$wsAlunni.Cells.Item(24, 3).NumberFormat = "#,##0"
$wsAlunni.Cells.Item(24, 4).Value = 3034501
$wsAlunni.Cells.Item(24, 4).NumberFormat = "#,##0"
$wsAlunni.Cells.Item(24, 5).Value = 123428
$wsAlunni.Cells.Item(24, 5).NumberFormat = "#,##0"
$wsAlunni.Cells.Item(24, 6).Value = 0.039
$wsAlunni.Cells.Item(24, 6).NumberFormat = "0.0%"

# ---------------------------------------------------------------------------
# Sheet "Personale scolastico" (sheet4.xml) -- add row 7
# ---------------------------------------------------------------------------
$wsPersonale = $wb.Worksheets.Item("Personale scolastico")

$wsPersonale.Cells.Item(7, 1).Value = "14-19 febbraio 2022"
$wsPersonale.Cells.Item(7, 2).Value = 775867
$wsPersonale.Cells.Item(7, 2).NumberFormat = "#,##0"
$wsPersonale.Cells.Item(7, 3).Value = 578258
$wsPersonale.Cells.Item(7, 3).NumberFormat = "#,##0"
$wsPersonale.Cells.Item(7, 4).Value = 0.745
$wsPersonale.Cells.Item(7, 4).NumberFormat = "0.0%"
$wsPersonale.Cells.Item(7, 5).Value = 557629
$wsPersonale.Cells.Item(7, 5).NumberFormat = "#,##0"
$wsPersonale.Cells.Item(7, 6).Value = 0.964
$wsPersonale.Cells.Item(7, 6).NumberFormat = "0.0%"
$wsPersonale.Cells.Item(7, 7).Value = 204526
$wsPersonale.Cells.Item(7, 7).NumberFormat = "#,##0"
$wsPersonale.Cells.Item(7, 8).Value = 153254
$wsPersonale.Cells.Item(7, 8).NumberFormat = "#,##0"
$wsPersonale.Cells.Item(7, 9).Value = 0.749
$wsPersonale.Cells.Item(7, 9).NumberFormat = "0.0%"
$wsPersonale.Cells.Item(7, 10).Value = 149129
$wsPersonale.Cells.Item(7, 10).NumberFormat = "#,##0"
$wsPersonale.Cells.Item(7, 11).Value = 0.973
$wsPersonale.Cells.Item(7, 11).NumberFormat = "0.0%"

# ---------------------------------------------------------------------------
# Selections per sheet (mirrors the final author view state)
# ---------------------------------------------------------------------------
$wsClassi.Range("A7").Select()
$wsAlunniPresenza.Range("F8").Select()
$wsAlunni.Range("C24").Select()

# Activate "Personale scolastico" last so it becomes the active/selected tab
# (workbook activeTab=3, sheet tabSelected moves here) and set its selection.
$wsPersonale.Activate()
$wsPersonale.Range("I8").Select()
